# Daily attendance processing - 2026-01-25 09:34:59
# Reverse the order of the comma-separated "Recorded By" entries in column G
# (e.g. "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if ($null -ne $val -and $val -ne "") {
        $parts = $val -split ", "
        if ($parts.Count -gt 1) {
            $reversed = $parts[($parts.Count - 1)..0]
            $cell.Value2 = [string]::Join(", ", $reversed)
        }
    }
}
